$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the whole "Make sure that ghost transition properly between..."
#    list-paragraph (its text, its _GoBack bookmark, and the paragraph mark).
# ---------------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute(
    "Make sure that ghost transition properly between Ghost states and ghost frightened modes.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$deadPara = $hit.Paragraphs(1)
$deadPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Turn the (now first) empty paragraph that follows "Scatter for 5
#    seconds, then switch to Chase mode permanently." into the new score
#    note, and give it its own _GoBack bookmark positioned right after the
#    text (mirroring the target markup:
#      <w:r><w:t>...</w:t></w:r><w:bookmarkStart .../><w:bookmarkEnd .../>).
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute(
    "Scatter for 5 seconds, then switch to Chase mode permanently.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $anchor.Paragraphs(1).Next()

# Write the text with a trailing sentinel so Find can locate the exact
# insertion point afterwards (a collapsed range placed exactly at the
# paragraph-mark boundary resolves ambiguously, so we anchor on a real
# character and strip it once the bookmark is placed).
$targetPara.Range.Text = "Ghost points when eaten: 200, 400, 800, 1600~"

$sentinel = $d.Content
$sentinel.Find.Execute("~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentinel.Collapse(1)
$d.Bookmarks.Add("_GoBack", $sentinel)

$cleanup = $d.Content
$cleanup.Find.Execute("~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cleanup.Text = ""
